$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item(1).Name = "GNG_TO-16511686821808066"
$wb.Worksheets.Item(2).Name = "NB_TO-16511686836356843"
$wb.Worksheets.Item(3).Name = "RS_TO-1651168683636682"
$wb.Worksheets.Item(4).Name = "TOL_TO-1651168683684683"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16511686837477736"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511686821428072.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168682163811.csv"
$ws1.Range("B4").Value = "go_stims-1651168682164806.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686821798425.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16511686824895203.csv"
$ws2.Range("B3").Value = "TB-16511686836177084.csv"
$ws2.Range("B4").Value = "ZB-match_0-1651168682241844.csv"
$ws2.Range("B5").Value = "TB-16511686830314932.csv"
$ws2.Range("B6").Value = "ZB-match_9-16511686821888373.csv"
$ws2.Range("B7").Value = "TB-16511686828115177.csv"
$ws2.Range("B8").Value = "OB-16511686823378098.csv"
$ws2.Range("B9").Value = "OB-16511686825864902.csv"
$ws2.Range("B10").Value = "ZB-match_4-16511686822698073.csv"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511686836517143.csv"
$ws4.Range("B3").Value = "ZM_stims-16511686836376834.csv"
$ws4.Range("B4").Value = "MM_stims-16511686836676784.csv"
$ws4.Range("B5").Value = "ZM_stims-16511686836517143.csv"
$ws4.Range("B6").Value = "MM_stims-16511686836837149.csv"
$ws4.Range("B7").Value = "ZM_stims-16511686836676784.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511686837317436.csv"
$ws5.Range("B3").Value = "SAT_stims-1651168683699738.csv"
$ws5.Range("B4").Value = "SAT_stims-16511686836876774.csv"
$ws5.Range("B5").Value = "vSAT_stims-16511686837157736.csv"
